$wb = $excel.ActiveWorkbook

$wsTracking = $wb.Worksheets.Item("MXTMS")
$wsConnect  = $wb.Worksheets.Item("Connect")

# Helper scratch cell used to force text-typed (shared-string) values into
# cells, mirroring the original "digits-as-text" OrderID/PickupID/CourierID
# columns (D2/E2 on MXTMS, AG2 on Connect) without disturbing their
# existing (default) cell formatting/style.
function Set-TextValue {
    param($range, [string]$text)

    $scratch = $range.Worksheet.Range("ZZ9999")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

# MXTMS!D2 (OrderID) and MXTMS!E2 (PickupID) roll forward to the next
# tracking numbers.
Set-TextValue $wsTracking.Range("D2") "11191070"
Set-TextValue $wsTracking.Range("E2") "7392424"

# Connect!AG2 (PickupID) rolls forward to match.
Set-TextValue $wsConnect.Range("AG2") "15595151"
